$d = $word.ActiveDocument
$insertStart = $d.Content.End
$r = $d.Content
$r.Collapse(0)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Sinespaciado"/><w:tabs><w:tab w:val="left" w:pos="5670"/></w:tabs><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:color w:val="4472C4" w:themeColor="accent1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>HUJG-001</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:color w:val="4472C4" w:themeColor="accent1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Clasificación de juegos ingresado</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:br/><w:t>Como:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Ad</w:t></w:r><w:r><w:t>ministrado del sistema</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Quiero: </w:t></w:r><w:r><w:t>Quiero ordenar cada juego por su género para que su búsqueda sea más fácil</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Para:</w:t></w:r><w:r><w:t xml:space="preserve"> Para poder realizar una búsqueda de un juego según su género.</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:color w:val="4472C4" w:themeColor="accent1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>HUJG-002 Filtrado por género en catálogo de juegos</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:color w:val="4472C4" w:themeColor="accent1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Como:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>U</w:t></w:r><w:r><w:t>suario registrado</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/></w:rPr><w:t>Quiero</w:t></w:r><w:r><w:t>: poder filtrar los juegos disponibles por género</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Para</w:t></w:r><w:r><w:t>: encontrar fácilmente títulos que se ajusten a mis preferencias personales</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:color w:val="4472C4" w:themeColor="accent1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>HUJG-00</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:color w:val="4472C4" w:themeColor="accent1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:color w:val="4472C4" w:themeColor="accent1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:color w:val="4472C4" w:themeColor="accent1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Estadísticas por género de juegos</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:color w:val="4472C4" w:themeColor="accent1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Como:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>A</w:t></w:r><w:r><w:t>nalista de datos del sistema</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/></w:rPr><w:t>Quiero</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t xml:space="preserve">Quiere </w:t></w:r><w:r><w:t>obtener estadísticas detalladas de juegos agrupados por género</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Para</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t xml:space="preserve">Para </w:t></w:r><w:r><w:t>analizar qué tipos de juegos tienen mayor aceptación entre los usuarios</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Sinespaciado"/><w:tabs><w:tab w:val="left" w:pos="5670"/></w:tabs><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'
$r.InsertXML($xml)

$searchRange = $d.Range($insertStart, $d.Content.End)

$f1 = $searchRange.Duplicate
$f1.Find.ClearFormatting()
$f1.Find.Execute("Quiero: poder filtrar", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($f1.Find.Found) {
  $run1 = $d.Range($f1.Start, $f1.Start + 6)
  $run1.Style = "normaltextrun"
} else {
  Write-Output "NOT FOUND 1"
}

$f2 = $searchRange.Duplicate
$f2.Find.ClearFormatting()
$f2.Find.Execute("Quiero: Quiere", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($f2.Find.Found) {
  $run2 = $d.Range($f2.Start, $f2.Start + 6)
  $run2.Style = "normaltextrun"
} else {
  Write-Output "NOT FOUND 2"
}

Write-Output "done"
